$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Month" column header and two new month values
$ws.Range("B1").Value = "Month"
$ws.Range("B2").Value = "Helmikuu"
$ws.Range("B3").Value = "Maaliskuu"

# Adjust column widths: A gets narrower (no longer auto bestFit), B gets a
# width to fit the new month labels.
$ws.Columns.Item(1).ColumnWidth = 44.6640625
$ws.Columns.Item(2).ColumnWidth = 9.88671875

# Select B3, matching the saved selection state
$ws.Range("B3").Select() | Out-Null
